$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.579.24"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.616.73"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "511.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.58"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  -2.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.629.39"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("E10").Value = "  +4.27%  "
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.074.71"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.517.45"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.63"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.628.88"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.77"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.51"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.62"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.17"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0843"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.41%  "
$ws.Range("E29").Value = "  -2.61%  "
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.45"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.62"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.81"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.73%  "
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("E36").Value = "  -2.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.888"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.49"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.848"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.39"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.77"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "294.72"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.626"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.101"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.996"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0555"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.81%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.86"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.90"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.003.84"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.44%  "
